# Weekly Fruta/Hortaliza price update for Agricola del Norte S.A. de Arica - Kiwi.
# The underlying data rows (2:33) were reshuffled/refreshed with new values;
# this script rewrites the cells that differ from the original workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45034
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25600
$ws.Range("Q2").Value = "`$/bandeja 18 kilos"
$ws.Range("S2").Value = 1422
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = 45086
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25500
$ws.Range("S3").Value = 1417

$ws.Range("D4").Value = 45086
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("S4").Value = 1139

$ws.Range("D5").Value = 45002
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("Q5").Value = "`$/bandeja 18 kilos"
$ws.Range("S5").Value = 1361
$ws.Range("T5").Value = 18

$ws.Range("D6").Value = 44629
$ws.Range("L6").Value = "Segunda"
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("S6").Value = 972

$ws.Range("D7").Value = 44323
$ws.Range("M7").Value = 270
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("S7").Value = 1194

$ws.Range("D8").Value = 45148
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 23000
$ws.Range("P8").Value = 22500
$ws.Range("Q8").Value = "`$/bandeja 18 kilos"
$ws.Range("S8").Value = 1250
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44291
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 17500
$ws.Range("S9").Value = 972

$ws.Range("D10").Value = 44491
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("Q10").Value = "`$/bandeja 10 kilos"
$ws.Range("S10").Value = 1450
$ws.Range("T10").Value = 10

$ws.Range("D11").Value = 44307
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19500
$ws.Range("S11").Value = 1083

$ws.Range("D12").Value = 45163
$ws.Range("M12").Value = 270
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("Q12").Value = "`$/bandeja 18 kilos"
$ws.Range("S12").Value = 1083
$ws.Range("T12").Value = 18

$ws.Range("D13").Value = 45169
$ws.Range("M13").Value = 270
$ws.Range("N13").Value = 27000
$ws.Range("O13").Value = 28000
$ws.Range("P13").Value = 27500
$ws.Range("S13").Value = 1528

$ws.Range("D14").Value = 45169
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 26000
$ws.Range("O14").Value = 26000
$ws.Range("P14").Value = 26000
$ws.Range("S14").Value = 1444

$ws.Range("D15").Value = 44489
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 26000
$ws.Range("O15").Value = 27000
$ws.Range("P15").Value = 26500
$ws.Range("S15").Value = 1472

$ws.Range("D16").Value = 44656
$ws.Range("M16").Value = 270
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19500
$ws.Range("Q16").Value = "`$/bandeja 18 kilos"
$ws.Range("S16").Value = 1083
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 45134
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 350
$ws.Range("N17").Value = 21000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21429
$ws.Range("S17").Value = 1190

$ws.Range("D18").Value = 44706
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 9000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 9500
$ws.Range("Q18").Value = "`$/bandeja 10 kilos"
$ws.Range("S18").Value = 950
$ws.Range("T18").Value = 10

$ws.Range("D19").Value = 44789
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 250
$ws.Range("N19").Value = 19000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 19500
$ws.Range("S19").Value = 1083

$ws.Range("D20").Value = 44616
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 16000
$ws.Range("O20").Value = 17000
$ws.Range("P20").Value = 16500
$ws.Range("Q20").Value = "`$/caja 18 kilos granel"
$ws.Range("S20").Value = 917

$ws.Range("D21").Value = 45127
$ws.Range("M21").Value = 200

$ws.Range("D22").Value = 44784
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 19000
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 19500
$ws.Range("S22").Value = 1083

$ws.Range("D23").Value = 44602
$ws.Range("M23").Value = 270

$ws.Range("D24").Value = 45069
$ws.Range("K24").Value = "Sin especificar"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 370
$ws.Range("P24").Value = 19486
$ws.Range("R24").Value = "Región Metropolitana"

$ws.Range("D25").Value = 44487
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 15000
$ws.Range("P25").Value = 14500
$ws.Range("Q25").Value = "`$/bandeja 10 kilos"
$ws.Range("S25").Value = 1450
$ws.Range("T25").Value = 10

$ws.Range("D26").Value = 45043
$ws.Range("K26").Value = "Hayward"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 21000
$ws.Range("O26").Value = 22000
$ws.Range("P26").Value = 21500
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 1194

$ws.Range("D27").Value = 44418
$ws.Range("M27").Value = 240
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 11000
$ws.Range("P27").Value = 10500
$ws.Range("Q27").Value = "`$/bandeja 10 kilos"
$ws.Range("S27").Value = 1050
$ws.Range("T27").Value = 10

$ws.Range("D28").Value = 44991
$ws.Range("M28").Value = 250
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 25000
$ws.Range("P28").Value = 24500
$ws.Range("Q28").Value = "`$/bandeja 18 kilos"
$ws.Range("S28").Value = 1361
$ws.Range("T28").Value = 18

$ws.Range("D29").Value = 45107
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 320
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 21000
$ws.Range("P29").Value = 20500
$ws.Range("S29").Value = 1139

$ws.Range("D30").Value = 44819
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 17000
$ws.Range("O30").Value = 18000
$ws.Range("P30").Value = 17500
$ws.Range("Q30").Value = "`$/bandeja 10 kilos"
$ws.Range("S30").Value = 1750
$ws.Range("T30").Value = 10

$ws.Range("D31").Value = 44614
$ws.Range("M31").Value = 250
$ws.Range("N31").Value = 20000
$ws.Range("O31").Value = 21000
$ws.Range("P31").Value = 20500
$ws.Range("S31").Value = 1139

$ws.Range("D32").Value = 44673
$ws.Range("L32").Value = "Especial"
$ws.Range("M32").Value = 400
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 15000
$ws.Range("P32").Value = 14500
$ws.Range("Q32").Value = "`$/bandeja 10 kilos"
$ws.Range("S32").Value = 1450
$ws.Range("T32").Value = 10

$ws.Range("D33").Value = 44263
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 250
$ws.Range("N33").Value = 21000
$ws.Range("O33").Value = 22000
$ws.Range("P33").Value = 21500
$ws.Range("Q33").Value = "`$/caja 18 kilos"
$ws.Range("S33").Value = 1194
